# "astrolabe measurements.xlsx" - Add files via upload
#
# The re-uploaded workbook trims the trailing all-zero helper rows
# (52-57) from Sheet1 and removes the previously-empty column H,
# which shifts the "Standard Deviation" column (with its formulas)
# from column I into column H. The view/selection and the column
# width formatting follow the data to their new locations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H was empty/unused; deleting it shifts column I (the
# "Standard Deviation" values/formulas and its custom width) left
# into column H.
$ws.Columns.Item(8).Delete()

# Rows 52-57 never had source data - they only carried left-over
# formulas that evaluated to 0. Remove them completely.
$ws.Range("A52:H57").EntireRow.Delete()

# Reset the view: no more scrolled-down topLeftCell, and the
# selection now sits on the relocated header cell, H1.
[void]$ws.Range("H1").Select()
